# Add a new weekly snapshot sheet "20191110" after the last existing sheet
# "20191103", using it as a template (same layout/formulas), then fill in
# the new week's numbers, and make the new sheet the active tab.

$wb = $excel.ActiveWorkbook

$srcName = "20191103"
$newName = "20191110"

$src = $wb.Worksheets.Item($srcName)

# Copy the template sheet so formulas / number formats / shared-formula
# structure all carry over exactly; Copy places the clone right after the
# source sheet and makes it the active sheet/tab.
$src.Copy($null, $src)

$ws = $wb.Worksheets.Item($src.Index + 1)
$ws.Name = $newName

# --- Position-distribution table (rows 1-14) ---------------------------
# Column A labels + column D buckets are already correct via the copy;
# only the count (B), percentage (C) and the summary row 14 change.

$ws.Range("B1").Value = 284
$ws.Range("C1").Value = 0.05

$ws.Range("B2").Value = 133
$ws.Range("C2").Value = 0.02

$ws.Range("B3").Value = 167
$ws.Range("C3").Value = 0.03

$ws.Range("B4").Value = 184
$ws.Range("C4").Value = 0.03

$ws.Range("B5").Value = 227
$ws.Range("C5").Value = 0.04

$ws.Range("B6").Value = 210
$ws.Range("C6").Value = 0.04

$ws.Range("B7").Value = 297
$ws.Range("C7").Value = 0.05

$ws.Range("B8").Value = 281
$ws.Range("C8").Value = 0.05

$ws.Range("B9").Value = 280
$ws.Range("C9").Value = 0.05

$ws.Range("B10").Value = 388
$ws.Range("C10").Value = 0.07

$ws.Range("B11").Value = 528
$ws.Range("C11").Value = 0.1

$ws.Range("B12").Value = 545
$ws.Range("C12").Value = 0.1

$ws.Range("B13").Value = 324
$ws.Range("C13").Value = 0.06

$ws.Range("B14").Value = 1349
$ws.Range("C14").Value = 0.25

# --- Sentiment poll (rows 20-23) ---------------------------------------
$ws.Range("A20").Value = "看多"
$ws.Range("B20").Value = 1255
$ws.Range("C20").Value = 0.24

$ws.Range("A21").Value = "看空 (已选)"
$ws.Range("B21").Value = 1960
$ws.Range("C21").Value = 0.37

$ws.Range("A22").Value = "看平"
$ws.Range("B22").Value = 620
$ws.Range("C22").Value = 0.11

$ws.Range("A23").Value = "我是来给卫斯理打Call的"
$ws.Range("B23").Value = 1360
$ws.Range("C23").Value = 0.26

$ws.Activate()
$null = $ws.Range("C17").Select()
